$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.734.46'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.601.18'
$ws.Range("E3").Value = '  +0.21%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("E5").Value = '  +0.20%  '

$ws.Range("E6").Value = '  -0.16%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("E9").Value = '  -0.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.64'
$ws.Range("E10").Value = '  +0.67%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.823.47'
$ws.Range("E12").Value = '  +0.07%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.584.74'
$ws.Range("E13").Value = '  -0.14%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.07'
$ws.Range("E14").Value = '  +0.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.526'
$ws.Range("E15").Value = '  +0.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.07'
$ws.Range("E16").Value = '  -0.39%  '

$ws.Range("E17").Value = '  -3.32%  '

$ws.Range("E18").Value = '  +0.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.18'
$ws.Range("E19").Value = '  +0.51%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '209.02'
$ws.Range("E20").Value = '  -0.14%  '

$ws.Range("E21").Value = '  +0.80%  '

$ws.Range("E22").Value = '  -3.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.04'
$ws.Range("E23").Value = '  +1.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.79'
$ws.Range("E24").Value = '  +0.48%  '

$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("E27").Value = '  -0.87%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.35'
$ws.Range("E28").Value = '  +0.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0509'
$ws.Range("E29").Value = '  -1.97%  '

$ws.Range("E30").Value = '  -0.21%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.27'
$ws.Range("E31").Value = '  +0.58%  '

$ws.Range("E32").Value = '  +0.46%  '

$ws.Range("E33").Value = '  +20.74%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.279.09'
$ws.Range("E34").Value = '  -0.50%  '

$ws.Range("E35").Value = '  +1.69%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.50'
$ws.Range("E36").Value = '  +0.46%  '

$ws.Range("E37").Value = '  -4.38%  '

$ws.Range("E38").Value = '  -1.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.826'
$ws.Range("E39").Value = '  -0.11%  '

$ws.Range("E40").Value = '  +0.58%  '

$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("E42").Value = '  -1.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '62.75'
$ws.Range("E43").Value = '  -0.75%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.735.73'
$ws.Range("E44").Value = '  +0.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.54'

$ws.Range("E46").Value = '  +0.28%  '

$ws.Range("E47").Value = '  +1.99%  '

$ws.Range("E48").Value = '  +1.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.53'
$ws.Range("E49").Value = '  +2.60%  '

$ws.Range("E50").Value = '  +0.14%  '

$ws.Range("E51").Value = '  +1.54%  '
